$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Find-ParaIndex($doc, $needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Change 1: the "graph LR" mermaid Source Code block under "项目总体" becomes
# a single First Paragraph styled paragraph reading "项目总体图 ".
# ---------------------------------------------------------------------------
$idx1 = Find-ParaIndex $d "graph LR"
$p1 = $d.Paragraphs.Item($idx1)
$xml1 = $pkgOpen + '<w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">项目总体图</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' + $pkgClose
$p1.Range.InsertXML($xml1) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: append "服务端模块图" to the paragraph that ends with
# "服务端使用RESTful HTTP接口", drop the whole "classDiagram" Source Code
# block that used to follow it, and turn the next paragraph ("使用接口：...")
# from First Paragraph style into Body Text style.
# ---------------------------------------------------------------------------
$idx2 = Find-ParaIndex $d "服务端使用RESTful HTTP接口"
$p2 = $d.Paragraphs.Item($idx2)
$xml2 = $pkgOpen + '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">模块及接口：</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">服务端使用RESTful HTTP接口</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">服务端模块图</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' + $pkgClose
$p2.Range.InsertXML($xml2) | Out-Null

$idx3 = Find-ParaIndex $d "classDiagram"
$p3 = $d.Paragraphs.Item($idx3)
$p3.Range.Delete() | Out-Null

$idx4 = Find-ParaIndex $d "使用接口："
$p4 = $d.Paragraphs.Item($idx4)
$xml4 = $pkgOpen + '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">使用接口：</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">+ 数据库</w:t></w:r></w:p>' + $pkgClose
$p4.Range.InsertXML($xml4) | Out-Null
